$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Add the new worksheet "coisas a fazer" right after Plan1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "coisas a fazer"

# 2. Fill in content - order matches how the shared strings ended up being appended
$ws2.Range("B1").Value = "Coisas que ainda faltam serem feitas no projeto"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Terminar a página de Inicio (História da Igreja e da Pastoral, e pensar e projetar como será a página)"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Colocar a Quant. De crismandos no card de grupos da crisma (pagina Grupos)"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "Criar página de perfi, com as informações do usuário (nele terá a opção de editar infor., foto de perfil (?), senha e email)"

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = "Página p/ relatar bugs"

$ws2.Range("A7").Value = 6
$ws2.Range("B7").Value = "Mehorar o digitamento das informações (por exemplo: permitir onde é para aceitar numero, só aceite numero e já fique formatado)"

$ws2.Range("A8").Value = 7
$ws2.Range("B8").Value = "Modificar algumas coisas no banco de dados (fazer melhoramentos e adicionar colunas ( adicionar colunas em ""crismandos"" como por exemplo estado civil, se possui filhos, cidade, etc))"

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "Melhorar ou modificar o sistema de login e cadastro de catequista"

# 3. Column width / wrap text / row heights
$ws2.Columns.Item(2).ColumnWidth = 50.8
$ws2.Range("B2:B8").WrapText = $true
$ws2.Range("B9:B13").WrapText = $true

$ws2.Rows.Item(1).RowHeight = 18
$ws2.Rows.Item(2).RowHeight = 28.8
$ws2.Rows.Item(3).RowHeight = 28.8
$ws2.Rows.Item(4).RowHeight = 28.8
$ws2.Rows.Item(5).RowHeight = 28.8
$ws2.Rows.Item(7).RowHeight = 43.2
$ws2.Rows.Item(8).RowHeight = 57.6

# 4. Selections (must set sheet1's selection first, sheet2's last so sheet2 ends active)
$ws1.Range("C3").Select() | Out-Null
$ws2.Range("B4").Select() | Out-Null
